$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 - Sep 28 2020
$ws.Range("D11").Value = 1454
$ws.Range("E11").Value = 722
$ws.Range("F11").Value = 2176
$ws.Range("G11").Value = 10
$ws.Range("E11").Style = "Bad"

# Row 12 - Sep 29 2020
$ws.Range("D12").Value = 2176
$ws.Range("E12").Value = -516
$ws.Range("F12").Value = 1660
$ws.Range("G12").Value = 13
$ws.Range("E12").Style = "Good"

# Update selection to J19
$ws.Range("J19").Select()
